$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2023" column (K), mirroring the formatting
# already used by the preceding "2022" column (J).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 581.79999999999995
$ws.Range("K5").Value = 571.20000000000005
$ws.Range("K6").Value = 584.1
